# "What is PSO? (intro ans inspiration added)"
#
# Slide 2 ("What is PSO?") - Content Placeholder 2: replace the single
# empty paragraph with the full bullet outline (intro + inspiration).
#
# Slide 4 ("Theoretical Analysis (...)") - Title: merge the two runs
# "Theoretical Analysis (" + "time complexity)" into a single run.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 2 : Content Placeholder 2
# ---------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$shape2 = $slide2.Shapes.Item("Content Placeholder 2")
$tf2 = $shape2.TextFrame
$tr2 = $tf2.TextRange

$CR = [char]13

# Paragraph plan: IndentLevel (PowerPoint 1-based; 2 => <a:pPr lvl="1"/>,
# 3 => <a:pPr lvl="2"/>, $null => leave default / no pPr) and its runs
# (text, font size in points or $null to leave default).
$paragraphs = @(
    @{ Level = 2; Runs = @(
        @{ Text = "A population based "; Size = $null },
        @{ Text = "stochastic "; Size = $null },
        @{ Text = "algorithm "; Size = $null },
        @{ Text = "based "; Size = $null },
        @{ Text = "on "; Size = $null },
        @{ Text = "metaheuristic"; Size = $null },
        @{ Text = " "; Size = $null },
        @{ Text = "approach."; Size = $null }
    )},
    @{ Level = 2; Runs = @(
        @{ Text = "Used in numerical "; Size = $null },
        @{ Text = "optimization "; Size = $null },
        @{ Text = "problems."; Size = $null }
    )},
    @{ Level = 2; Runs = @(
        @{ Text = "An inspiration from "; Size = $null },
        @{ Text = "social foraging behaviors of "; Size = $null },
        @{ Text = "animals.."; Size = $null }
    )},
    @{ Level = 2; Runs = @(
        @{ Text = "A"; Size = $null },
        @{ Text = "lso "; Size = $null },
        @{ Text = "classified as swarm intelligence algorithm "; Size = $null },
        @{ Text = "like,"; Size = $null }
    )},
    @{ Level = 3; Runs = @(
        @{ Text = "bacterial "; Size = 20 },
        @{ Text = "foraging "; Size = 20 },
        @{ Text = "algorithm"; Size = 20 }
    )},
    @{ Level = 3; Runs = @(
        @{ Text = "ant "; Size = 20 },
        @{ Text = "colony algorithm etc."; Size = 20 }
    )},
    @{ Level = 2; Runs = @(
        @{ Text = "Example"; Size = $null },
        @{ Text = ": predicting score of a football team using a math equation."; Size = $null }
    )},
    @{ Level = 2; Runs = @() },
    @{ Level = $null; Runs = @() }
)

# Build the full text blob (paragraphs separated by CR) in one shot so the
# shape materialises every paragraph at once.
$full = ""
for ($i = 0; $i -lt $paragraphs.Count; $i++) {
    $para = $paragraphs[$i]
    foreach ($run in $para.Runs) {
        $full += $run.Text
    }
    if ($i -lt ($paragraphs.Count - 1)) {
        $full += $CR
    }
}
$tr2.Text = $full

# Force PowerPoint-style "shrink text on overflow" autofit, matching
# <a:bodyPr><a:normAutofit/></a:bodyPr> in the target deck.
$tf2.AutoSize = 2

# Walk the paragraphs again, this time applying indent level + per-run
# font (Times New Roman / size) using precise character offsets.
$pos = 1
for ($i = 0; $i -lt $paragraphs.Count; $i++) {
    $para = $paragraphs[$i]
    $paraStart = $pos
    $paraLen = 0
    foreach ($run in $para.Runs) {
        $paraLen += $run.Text.Length
    }

    if ($paraLen -gt 0) {
        foreach ($run in $para.Runs) {
            $len = $run.Text.Length
            $rng = $tr2.Characters($pos, $len)
            $rng.Font.Name = "Times New Roman"
            $rng.Font.NameComplexScript = "Times New Roman"
            if ($run.Size -ne $null) {
                $rng.Font.Size = $run.Size
            }
            $pos += $len
        }
    } else {
        # Empty paragraph: still tag the (zero-length) run with the
        # correct trailing-run font so the end-of-paragraph formatting
        # matches (last paragraph uses sz=20).
        $rng = $tr2.Characters($pos, 0)
        $rng.Font.Name = "Times New Roman"
        $rng.Font.NameComplexScript = "Times New Roman"
        if ($i -eq ($paragraphs.Count - 1)) {
            $rng.Font.Size = 20
        }
    }

    if ($para.Level -ne $null) {
        $paraRange = $tr2.Paragraphs($i + 1)
        $paraRange.IndentLevel = $para.Level
    }

    # Skip past the paragraph break.
    $pos = $paraStart + $paraLen + 1
}

# ---------------------------------------------------------------------
# Slide 4 : Title - merge "Theoretical Analysis (" + "time complexity)"
# ---------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$shape4 = $slide4.Shapes.Item("Title 1")
$tr4 = $shape4.TextFrame.TextRange

$tailLen = "time complexity)".Length
$headLen = $tr4.Length - $tailLen

$tail = $tr4.Characters($headLen + 1, $tailLen)
$tail.Text = ""

$head = $tr4.Characters(1, $headLen)
$head.InsertAfter("time complexity)") | Out-Null
